$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.080435
$ws.Range("H2").Value = 24.241305
$ws.Range("I2").Value = 0.1496988574979475
$ws.Range("J2").Value = 0.1496988574979476
$ws.Range("M2").Value = 0.6574793333333333
$ws.Range("N2").Value = 1.972438
$ws.Range("O2").Value = 0.04234443143670402
$ws.Range("P2").Value = 0.04234443143670403
$ws.Range("Q2").Value = 5.312719016843332
$ws.Range("R2").Value = 47.81447115159
$ws.Range("S2").Value = 0.006338913007474765
$ws.Range("T2").Value = 0.006338913007474767

$ws.Range("G3").Value = 8.080435
$ws.Range("H3").Value = 24.241305
$ws.Range("I3").Value = 0.1496988574979475
$ws.Range("J3").Value = 0.1496988574979476
$ws.Range("O3").Value = 0.1192373589365509
$ws.Range("P3").Value = 0.119237358936551
$ws.Range("Q3").Value = 14.960044635085
$ws.Range("R3").Value = 134.640401715765
$ws.Range("S3").Value = 0.01784969640387436
$ws.Range("T3").Value = 0.01784969640387437

$ws.Range("G4").Value = 8.080435
$ws.Range("H4").Value = 24.241305
$ws.Range("I4").Value = 0.1496988574979475
$ws.Range("J4").Value = 0.1496988574979476
$ws.Range("M4").Value = 5.370269333333333
$ws.Range("N4").Value = 16.110808
$ws.Range("O4").Value = 0.3458679080132824
$ws.Range("P4").Value = 0.3458679080132824
$ws.Range("Q4").Value = 43.39411228049332
$ws.Range("R4").Value = 390.54701052444
$ws.Range("S4").Value = 0.05177603067479358
$ws.Range("T4").Value = 0.0517760306747936

$ws.Range("G5").Value = 8.080435
$ws.Range("H5").Value = 24.241305
$ws.Range("I5").Value = 0.1496988574979475
$ws.Range("J5").Value = 0.1496988574979476
$ws.Range("M5").Value = 1.801189666666667
$ws.Range("N5").Value = 5.403569
$ws.Range("O5").Value = 0.1160041821512257
$ws.Range("P5").Value = 0.1160041821512257
$ws.Range("Q5").Value = 14.55439602417167
$ws.Range("R5").Value = 130.989564217545
$ws.Range("S5").Value = 0.01736569353302229
$ws.Range("T5").Value = 0.01736569353302229

$ws.Range("G6").Value = 8.080435
$ws.Range("H6").Value = 24.241305
$ws.Range("I6").Value = 0.1496988574979475
$ws.Range("J6").Value = 0.1496988574979476
$ws.Range("M6").Value = 5.846608
$ws.Range("N6").Value = 17.539824
$ws.Range("O6").Value = 0.3765461194622369
$ws.Range("P6").Value = 0.376546119462237
$ws.Range("Q6").Value = 47.24313591447999
$ws.Range("R6").Value = 425.18822323032
$ws.Range("S6").Value = 0.05636852387878254
$ws.Range("T6").Value = 0.05636852387878255

$ws.Range("I7").Value = 0.2404784903431001
$ws.Range("J7").Value = 0.2404784903431001
$ws.Range("M7").Value = 0.6574793333333333
$ws.Range("N7").Value = 1.972438
$ws.Range("O7").Value = 0.04234443143670402
$ws.Range("P7").Value = 0.04234443143670403
$ws.Range("Q7").Value = 8.534431525672002
$ws.Range("R7").Value = 76.809883731048
$ws.Range("S7").Value = 0.01018292494633549
$ws.Range("T7").Value = 0.0101829249463355

$ws.Range("I8").Value = 0.2404784903431001
$ws.Range("J8").Value = 0.2404784903431001
$ws.Range("O8").Value = 0.1192373589365509
$ws.Range("P8").Value = 0.119237358936551
$ws.Range("Q8").Value = 24.032040120012
$ws.Range("S8").Value = 0.02867402006956013
$ws.Range("T8").Value = 0.02867402006956014

$ws.Range("I9").Value = 0.2404784903431001
$ws.Range("J9").Value = 0.2404784903431001
$ws.Range("M9").Value = 5.370269333333333
$ws.Range("N9").Value = 16.110808
$ws.Range("O9").Value = 0.3458679080132824
$ws.Range("P9").Value = 0.3458679080132824
$ws.Range("Q9").Value = 69.70895292995201
$ws.Range("R9").Value = 627.380576369568
$ws.Range("S9").Value = 0.08317379237716038
$ws.Range("T9").Value = 0.08317379237716038

$ws.Range("I10").Value = 0.2404784903431001
$ws.Range("J10").Value = 0.2404784903431001
$ws.Range("M10").Value = 1.801189666666667
$ws.Range("N10").Value = 5.403569
$ws.Range("O10").Value = 0.1160041821512257
$ws.Range("P10").Value = 0.1160041821512257
$ws.Range("Q10").Value = 23.380400106236
$ws.Range("R10").Value = 210.423600956124
$ws.Range("S10").Value = 0.02789651059721276
$ws.Range("T10").Value = 0.02789651059721276

$ws.Range("I11").Value = 0.2404784903431001
$ws.Range("J11").Value = 0.2404784903431001
$ws.Range("M11").Value = 5.846608
$ws.Range("N11").Value = 17.539824
$ws.Range("O11").Value = 0.3765461194622369
$ws.Range("P11").Value = 0.376546119462237
$ws.Range("Q11").Value = 75.89208223545602
$ws.Range("R11").Value = 683.0287401191041
$ws.Range("S11").Value = 0.09055124235283138
$ws.Range("T11").Value = 0.09055124235283139

$ws.Range("G12").Value = 15.25749233333333
$ws.Range("H12").Value = 45.772477
$ws.Range("I12").Value = 0.2826616599952471
$ws.Range("J12").Value = 0.2826616599952471
$ws.Range("M12").Value = 0.6574793333333333
$ws.Range("N12").Value = 1.972438
$ws.Range("O12").Value = 0.04234443143670402
$ws.Range("P12").Value = 0.04234443143670403
$ws.Range("Q12").Value = 10.03148588765844
$ws.Range("R12").Value = 90.283372988926
$ws.Range("S12").Value = 0.01196914728145368
$ws.Range("T12").Value = 0.01196914728145368

$ws.Range("G13").Value = 15.25749233333333
$ws.Range("H13").Value = 45.772477
$ws.Range("I13").Value = 0.2826616599952471
$ws.Range("J13").Value = 0.2826616599952471
$ws.Range("O13").Value = 0.1192373589365509
$ws.Range("P13").Value = 0.119237358936551
$ws.Range("Q13").Value = 28.24758398850233
$ws.Range("R13").Value = 254.228255896521
$ws.Range("S13").Value = 0.0337038298104546
$ws.Range("T13").Value = 0.0337038298104546

$ws.Range("G14").Value = 15.25749233333333
$ws.Range("H14").Value = 45.772477
$ws.Range("I14").Value = 0.2826616599952471
$ws.Range("J14").Value = 0.2826616599952471
$ws.Range("M14").Value = 5.370269333333333
$ws.Range("N14").Value = 16.110808
$ws.Range("O14").Value = 0.3458679080132824
$ws.Range("P14").Value = 0.3458679080132824
$ws.Range("Q14").Value = 81.93684318126844
$ws.Range("R14").Value = 737.4315886314159
$ws.Range("S14").Value = 0.09776359701811781
$ws.Range("T14").Value = 0.09776359701811781

$ws.Range("G15").Value = 15.25749233333333
$ws.Range("H15").Value = 45.772477
$ws.Range("I15").Value = 0.2826616599952471
$ws.Range("J15").Value = 0.2826616599952471
$ws.Range("M15").Value = 1.801189666666667
$ws.Range("N15").Value = 5.403569
$ws.Range("O15").Value = 0.1160041821512257
$ws.Range("P15").Value = 0.1160041821512257
$ws.Range("Q15").Value = 27.48163753004589
$ws.Range("R15").Value = 247.334737770413
$ws.Range("S15").Value = 0.03278993469325647
$ws.Range("T15").Value = 0.03278993469325647

$ws.Range("G16").Value = 15.25749233333333
$ws.Range("H16").Value = 45.772477
$ws.Range("I16").Value = 0.2826616599952471
$ws.Range("J16").Value = 0.2826616599952471
$ws.Range("M16").Value = 5.846608
$ws.Range("N16").Value = 17.539824
$ws.Range("O16").Value = 0.3765461194622369
$ws.Range("P16").Value = 0.376546119462237
$ws.Range("Q16").Value = 89.20457673600534
$ws.Range("R16").Value = 802.841190624048
$ws.Range("S16").Value = 0.1064351511919645
$ws.Range("T16").Value = 0.1064351511919645

$ws.Range("G17").Value = 4.142925
$ws.Range("H17").Value = 12.428775
$ws.Range("I17").Value = 0.07675219702895753
$ws.Range("J17").Value = 0.07675219702895753
$ws.Range("M17").Value = 0.6574793333333333
$ws.Range("N17").Value = 1.972438
$ws.Range("O17").Value = 0.04234443143670402
$ws.Range("P17").Value = 0.04234443143670403
$ws.Range("Q17").Value = 2.72388756705
$ws.Range("R17").Value = 24.51498810345
$ws.Range("S17").Value = 0.00325002814470909
$ws.Range("T17").Value = 0.00325002814470909

$ws.Range("G18").Value = 4.142925
$ws.Range("H18").Value = 12.428775
$ws.Range("I18").Value = 0.07675219702895753
$ws.Range("J18").Value = 0.07675219702895753
$ws.Range("O18").Value = 0.1192373589365509
$ws.Range("P18").Value = 0.119237358936551
$ws.Range("Q18").Value = 7.670174058674999
$ws.Range("R18").Value = 69.03156652807499
$ws.Range("S18").Value = 0.009151729266310688
$ws.Range("T18").Value = 0.00915172926631069

$ws.Range("G19").Value = 4.142925
$ws.Range("H19").Value = 12.428775
$ws.Range("I19").Value = 0.07675219702895753
$ws.Range("J19").Value = 0.07675219702895753
$ws.Range("M19").Value = 5.370269333333333
$ws.Range("N19").Value = 16.110808
$ws.Range("O19").Value = 0.3458679080132824
$ws.Range("P19").Value = 0.3458679080132824
$ws.Range("Q19").Value = 22.2486230778
$ws.Range("R19").Value = 200.2376077002
$ws.Range("S19").Value = 0.02654612182182881
$ws.Range("T19").Value = 0.02654612182182881

$ws.Range("G20").Value = 4.142925
$ws.Range("H20").Value = 12.428775
$ws.Range("I20").Value = 0.07675219702895753
$ws.Range("J20").Value = 0.07675219702895753
$ws.Range("M20").Value = 1.801189666666667
$ws.Range("N20").Value = 5.403569
$ws.Range("O20").Value = 0.1160041821512257
$ws.Range("P20").Value = 0.1160041821512257
$ws.Range("Q20").Value = 7.462193699775
$ws.Range("R20").Value = 67.159743297975
$ws.Range("S20").Value = 0.008903575844653953
$ws.Range("T20").Value = 0.008903575844653955

$ws.Range("G21").Value = 4.142925
$ws.Range("H21").Value = 12.428775
$ws.Range("I21").Value = 0.07675219702895753
$ws.Range("J21").Value = 0.07675219702895753
$ws.Range("M21").Value = 5.846608
$ws.Range("N21").Value = 17.539824
$ws.Range("O21").Value = 0.3765461194622369
$ws.Range("P21").Value = 0.376546119462237
$ws.Range("Q21").Value = 24.2220584484
$ws.Range("R21").Value = 217.9985260356
$ws.Range("S21").Value = 0.02890074195145498
$ws.Range("T21").Value = 0.02890074195145499

$ws.Range("G22").Value = 13.51654933333334
$ws.Range("H22").Value = 40.549648
$ws.Range("I22").Value = 0.2504087951347477
$ws.Range("J22").Value = 0.2504087951347477
$ws.Range("M22").Value = 0.6574793333333333
$ws.Range("N22").Value = 1.972438
$ws.Range("O22").Value = 0.04234443143670402
$ws.Range("P22").Value = 0.04234443143670403
$ws.Range("Q22").Value = 8.886851844647111
$ws.Range("R22").Value = 79.981666601824
$ws.Range("S22").Value = 0.01060341805673099
$ws.Range("T22").Value = 0.01060341805673099

$ws.Range("G23").Value = 13.51654933333334
$ws.Range("H23").Value = 40.549648
$ws.Range("I23").Value = 0.2504087951347477
$ws.Range("J23").Value = 0.2504087951347477
$ws.Range("O23").Value = 0.1192373589365509
$ws.Range("P23").Value = 0.119237358936551
$ws.Range("Q23").Value = 25.02441778678934
$ws.Range("R23").Value = 225.219760081104
$ws.Range("S23").Value = 0.02985808338635116
$ws.Range("T23").Value = 0.02985808338635116

$ws.Range("G24").Value = 13.51654933333334
$ws.Range("H24").Value = 40.549648
$ws.Range("I24").Value = 0.2504087951347477
$ws.Range("J24").Value = 0.2504087951347477
$ws.Range("M24").Value = 5.370269333333333
$ws.Range("N24").Value = 16.110808
$ws.Range("O24").Value = 0.3458679080132824
$ws.Range("P24").Value = 0.3458679080132824
$ws.Range("Q24").Value = 72.58751037728712
$ws.Range("R24").Value = 653.2875933955841
$ws.Range("S24").Value = 0.08660836612138179
$ws.Range("T24").Value = 0.08660836612138179

$ws.Range("G25").Value = 13.51654933333334
$ws.Range("H25").Value = 40.549648
$ws.Range("I25").Value = 0.2504087951347477
$ws.Range("J25").Value = 0.2504087951347477
$ws.Range("M25").Value = 1.801189666666667
$ws.Range("N25").Value = 5.403569
$ws.Range("O25").Value = 0.1160041821512257
$ws.Range("P25").Value = 0.1160041821512257
$ws.Range("Q25").Value = 24.34586898819023
$ws.Range("R25").Value = 219.112820893712
$ws.Range("S25").Value = 0.02904846748308023
$ws.Range("T25").Value = 0.02904846748308023

$ws.Range("G26").Value = 13.51654933333334
$ws.Range("H26").Value = 40.549648
$ws.Range("I26").Value = 0.2504087951347477
$ws.Range("J26").Value = 0.2504087951347477
$ws.Range("M26").Value = 5.846608
$ws.Range("N26").Value = 17.539824
$ws.Range("O26").Value = 0.3765461194622369
$ws.Range("P26").Value = 0.376546119462237
$ws.Range("Q26").Value = 79.02596546466134
$ws.Range("R26").Value = 711.2336891819521
$ws.Range("S26").Value = 0.09429046008720351
$ws.Range("T26").Value = 0.09429046008720353
